$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number would be auto-converted by Excel
# from text to a numeric type. The source data keeps these as text (e.g. "215.55"),
# so force text storage via NumberFormat, then restore the default style so no
# stray formatting is left behind.

$ws.Range("D2").Value = '25.802.36'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.638.01'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -0.84%  '
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0792'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '1.862.66'
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("D14").Value = '1.637.56'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '25.823.02'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '193.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +3.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.64%  '
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("E31").Value = '  -1.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.908'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = '1.134.92'
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E38").Value = '  -2.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.544'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.63%  '
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.45%  '
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").Value = '1.771.97'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").Value = '0.0₆0113'
$ws.Range("E46").Value = '  +2.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("E48").Value = '  -1.54%  '
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.42'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.69%  '
